$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$tcs = $nm.Theme.ThemeColorScheme
$c1 = $tcs.Item(1)
$c1.RGB = 255
Write-Output ("after RGB=" + $tcs.Item(1).RGB)
